# Updated schedule, added TAs
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Move the room-map hyperlink text for the 17/11 block from K3 up to K2,
#     matching the other day blocks where the link sits on the date row,
#     and copy K3's formatting onto K2 so the style (s="1") moves with it.
$ws.Range("K3").Copy() | Out-Null
$ws.Range("K2").PasteSpecial(-4122) | Out-Null
$ws.Range("K2").Value = $ws.Range("K3").Value2
$ws.Range("K3").Clear()

# --- Add teaching-assistant (TA) initials in the new "assistant" column (G)
$ws.Range("G4").Value = "KG"
$ws.Range("G7").Value = "KG,SS"
$ws.Range("G9").Value = "JH,KG"
$ws.Range("G12").Value = "GD,JH"
$ws.Range("G14").Value = "GD,KG"
$ws.Range("G18").Value = "AJ,SD"
$ws.Range("G20").Value = "JH,SD"
$ws.Range("G23").Value = "GD,SS"
$ws.Range("G24").Value = "GD,SS"

# --- Replace instructor "Stephan Nylander" with "Elin Kronander" on 21/11,
#     and add her to the NBIS closing slot as well
$ws.Range("F26").Value = "Elin Kronander"
$ws.Range("F27").Value = "Elin Kronander"

# --- Update the active selection left behind in the sheet view
$ws.Range("G30").Select() | Out-Null
